$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(221).Insert()
$ws.Cells.Item(221, 1).Value = 7
$ws.Cells.Item(221, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(221, 3).Value = "Ñuble"
$ws.Cells.Item(221, 4).Value = 44841
$ws.Cells.Item(221, 5).Value = 16
$ws.Cells.Item(221, 6).Value = 100112009
$ws.Cells.Item(221, 7).Value = "Acelga"
$ws.Cells.Item(221, 8).Value = "Sin especificar"
$ws.Cells.Item(221, 9).Value = "Segunda"
$ws.Cells.Item(221, 10).Value = 200
$ws.Cells.Item(221, 11).Value = 600
$ws.Cells.Item(221, 12).Value = 600
$ws.Cells.Item(221, 13).Value = 600
$ws.Cells.Item(221, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(221, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(221, 16).Value = 600
$ws.Cells.Item(221, 17).Value = 1
$ws.Cells.Item(221, 18).Value = "Hortaliza"
